# Update scripts with new TPM values.
#
# The sheet originally held 5 "sending cluster" rows (ECs, FAPs,
# Inflammatory-Mac, MuSCs, Neutrophils) all signalling Ptn -> Alk into the
# FAPs target cluster. The refreshed TPM run keeps only the FAPs and MuSCs
# sending-cluster rows (with recomputed metrics) and drops the other three
# rows entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three obsolete sending-cluster rows (ECs' row 4 data that used to
# belong to Inflammatory-Mac/MuSCs/Neutrophils). Rows shift up as each is
# removed, so the same row index is deleted three times.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# Row 2 (was ECs -> FAPs/Alk/FAPs) becomes the FAPs sending-cluster row with
# refreshed TPM-derived metrics.
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ptn"
$ws.Range("C2").Value = "Alk"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.807741666666667
$ws.Range("H2").Value = 8.423225
$ws.Range("I2").Value = 0.3451213243389217
$ws.Range("J2").Value = 0.4414969759863614
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03425633333333333
$ws.Range("N2").Value = 0.102769
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.09618293444722223
$ws.Range("R2").Value = 0.8656464100250001
$ws.Range("S2").Value = 0.3451213243389217
$ws.Range("T2").Value = 0.4414969759863614

# Row 3 (was FAPs -> FAPs/Alk/FAPs) becomes the MuSCs sending-cluster row
# with refreshed TPM-derived metrics.
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Ptn"
$ws.Range("C3").Value = "Alk"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.327779
$ws.Range("H3").Value = 10.655558
$ws.Range("I3").Value = 0.6548786756610784
$ws.Range("J3").Value = 0.5585030240136385
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03425633333333333
$ws.Range("N3").Value = 0.102769
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.1825101733503333
$ws.Range("R3").Value = 1.095061040102
$ws.Range("S3").Value = 0.6548786756610784
$ws.Range("T3").Value = 0.5585030240136385
